# Apply updated "dSF" (column F) values on the active worksheet.
# These represent a repull/recalculation of the dSF stat for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -7
    5  = -6
    7  = 6
    10 = 6
    11 = 0
    12 = 14
    13 = -6
    14 = -4
    15 = -3
    16 = -3
    19 = 4
    21 = -8
    22 = -1
    23 = -1
    24 = -1
    25 = 2
    26 = -2
    34 = 7
    35 = -2
    37 = 0
    38 = -2
    39 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
